$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 2360.32395
$ws.Range("C2").Value = 6151.921500000001
$ws.Range("D2").Value = 0.8878
$ws.Range("E2").Value = 4474.151599999999
$ws.Range("F2").Value = 630074.3825000001
$ws.Range("G2").Value = 553200.78605
$ws.Range("H2").Value = 76873.59650000001
$ws.Range("I2").Value = 33156.82405
$ws.Range("J2").Value = 43716.77245
$ws.Range("K2").Value = 799697.4858
$ws.Range("L2").Value = 553200.799
$ws.Range("M2").Value = 246496.6868
$ws.Range("N2").Value = 74864.5508
$ws.Range("O2").Value = 171632.1358
$ws.Range("B3").Value = 2475.847
$ws.Range("C3").Value = 6526.769
$ws.Range("E3").Value = 2776.521
$ws.Range("F3").Value = 670863.079
$ws.Range("G3").Value = 568305.581
$ws.Range("H3").Value = 102557.498
$ws.Range("I3").Value = 49350.841
$ws.Range("J3").Value = 53206.657
$ws.Range("K3").Value = 786076.3017999999
$ws.Range("L3").Value = 568305.598
$ws.Range("M3").Value = 217770.7038
$ws.Range("N3").Value = 80335.28540000001
$ws.Range("O3").Value = 137435.4182
$ws.Range("B4").Value = 4453.662
$ws.Range("C4").Value = 5543.459
$ws.Range("E4").Value = 4738.919
$ws.Range("F4").Value = 909377.991
$ws.Range("G4").Value = 766298.265
$ws.Range("H4").Value = 143079.726
$ws.Range("I4").Value = 57813.58
$ws.Range("J4").Value = 85266.14599999999
$ws.Range("K4").Value = 811179.3924
$ws.Range("L4").Value = 774147.678
$ws.Range("M4").Value = 37031.7144
$ws.Range("N4").Value = 20358.9328
$ws.Range("O4").Value = 16672.7818
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 2359.05995
$ws.Range("C2").Value = 6231.029950000002
$ws.Range("D2").Value = 0.8878
$ws.Range("E2").Value = 4491.621249999999
$ws.Range("F2").Value = 635986.0139
$ws.Range("G2").Value = 556239.7016499999
$ws.Range("H2").Value = 79746.31245
$ws.Range("I2").Value = 31680.05055
$ws.Range("J2").Value = 48066.2619
$ws.Range("K2").Value = 746784.1983999999
$ws.Range("L2").Value = 556239.7169999999
$ws.Range("M2").Value = 190544.4814
$ws.Range("N2").Value = 58208.5764
$ws.Range("O2").Value = 132335.905
$ws.Range("B3").Value = 2596.697
$ws.Range("C3").Value = 6412.722
$ws.Range("F3").Value = 681037.809
$ws.Range("G3").Value = 580646.203
$ws.Range("H3").Value = 100391.606
$ws.Range("I3").Value = 46157.504
$ws.Range("J3").Value = 54234.102
$ws.Range("K3").Value = 738858.0728
$ws.Range("L3").Value = 580646.159
$ws.Range("M3").Value = 158211.9138
$ws.Range("N3").Value = 58572.50919999999
$ws.Range("O3").Value = 99639.40479999999
$ws.Range("B4").Value = 3680.941
$ws.Range("C4").Value = 5624.402
$ws.Range("E4").Value = 3221.5
$ws.Range("F4").Value = 791841.416
$ws.Range("G4").Value = 668535.621
$ws.Range("H4").Value = 123305.795
$ws.Range("I4").Value = 31231.031
$ws.Range("J4").Value = 92074.764
$ws.Range("K4").Value = 755955.9550000001
$ws.Range("L4").Value = 676385.035
$ws.Range("M4").Value = 79570.92
$ws.Range("N4").Value = 33297.649
$ws.Range("O4").Value = 46273.271
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 2503.9123
$ws.Range("C2").Value = 6187.141500000001
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 4098.746800000001
$ws.Range("F2").Value = 653832.27975
$ws.Range("G2").Value = 568367.3552999999
$ws.Range("H2").Value = 85464.92449999999
$ws.Range("I2").Value = 35922.63370000001
$ws.Range("J2").Value = 49542.2908
$ws.Range("K2").Value = 659707.4276000001
$ws.Range("L2").Value = 568367.368
$ws.Range("M2").Value = 91340.05959999999
$ws.Range("N2").Value = 42926.9698
$ws.Range("O2").Value = 48413.09
$ws.Range("B3").Value = 2664.22
$ws.Range("C3").Value = 6388.932
$ws.Range("E3").Value = 3142.458
$ws.Range("F3").Value = 700933.33
$ws.Range("G3").Value = 587548.3689999999
$ws.Range("H3").Value = 113384.961
$ws.Range("I3").Value = 47515.456
$ws.Range("J3").Value = 65869.505
$ws.Range("K3").Value = 659128.7048000001
$ws.Range("L3").Value = 587548.432
$ws.Range("M3").Value = 71580.27280000001
$ws.Range("N3").Value = 41828.311
$ws.Range("O3").Value = 29751.9618
$ws.Range("B4").Value = 4453.662
$ws.Range("C4").Value = 5543.459
$ws.Range("E4").Value = 4738.919
$ws.Range("F4").Value = 909377.991
$ws.Range("G4").Value = 766298.265
$ws.Range("H4").Value = 143079.726
$ws.Range("I4").Value = 57813.58
$ws.Range("J4").Value = 85266.14599999999
$ws.Range("K4").Value = 781604.7429999999
$ws.Range("L4").Value = 774147.678
$ws.Range("M4").Value = 7457.065
$ws.Range("N4").Value = 6426.1798
$ws.Range("O4").Value = 1030.8852
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 2616.641
$ws.Range("C2").Value = 6369.8743
$ws.Range("D2").Value = 0.8878
$ws.Range("E2").Value = 3824.308
$ws.Range("F2").Value = 673255.26105
$ws.Range("G2").Value = 586425.8946999999
$ws.Range("H2").Value = 86829.36649999997
$ws.Range("I2").Value = 39486.67635
$ws.Range("J2").Value = 47342.69014999999
$ws.Range("K2").Value = 625794.1323999999
$ws.Range("L2").Value = 586425.901
$ws.Range("M2").Value = 39368.2314
$ws.Range("N2").Value = 16678.635
$ws.Range("O2").Value = 22689.5964
$ws.Range("B3").Value = 2691.734
$ws.Range("C3").Value = 6595.794
$ws.Range("E3").Value = 3056.044
$ws.Range("F3").Value = 708209.27
$ws.Range("G3").Value = 598073.144
$ws.Range("H3").Value = 110136.125
$ws.Range("I3").Value = 49926.072
$ws.Range("J3").Value = 60210.053
$ws.Range("K3").Value = 633550.0596
$ws.Range("L3").Value = 598073.154
$ws.Range("M3").Value = 35476.9056
$ws.Range("N3").Value = 16256.2612
$ws.Range("O3").Value = 19220.6444
$ws.Range("B4").Value = 4453.662
$ws.Range("C4").Value = 5543.459
$ws.Range("E4").Value = 4738.919
$ws.Range("F4").Value = 909377.991
$ws.Range("G4").Value = 766298.265
$ws.Range("H4").Value = 143079.726
$ws.Range("I4").Value = 57813.58
$ws.Range("J4").Value = 85266.14599999999
$ws.Range("K4").Value = 785329.444
$ws.Range("L4").Value = 774147.678
$ws.Range("M4").Value = 11181.766
$ws.Range("N4").Value = 4579.151599999999
$ws.Range("O4").Value = 6602.6144
$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = 2545.495
$ws.Range("C2").Value = 6380.50295
$ws.Range("D2").Value = 0.8878
$ws.Range("E2").Value = 3506.086349999999
$ws.Range("F2").Value = 666123.5179999999
$ws.Range("G2").Value = 576207.5757
$ws.Range("H2").Value = 89915.94245
$ws.Range("I2").Value = 41518.32995
$ws.Range("J2").Value = 48397.6125
$ws.Range("K2").Value = 666133.5338
$ws.Range("L2").Value = 576207.5820000001
$ws.Range("M2").Value = 89925.9518
$ws.Range("N2").Value = 50248.1024
$ws.Range("O2").Value = 39677.8492
$ws.Range("B3").Value = 2659.375
$ws.Range("C3").Value = 6435.508
$ws.Range("E3").Value = 3343.039
$ws.Range("F3").Value = 699250.143
$ws.Range("G3").Value = 590290.827
$ws.Range("H3").Value = 108959.316
$ws.Range("I3").Value = 44081.989
$ws.Range("J3").Value = 64877.326
$ws.Range("K3").Value = 665959.287
$ws.Range("L3").Value = 590290.847
$ws.Range("M3").Value = 75668.43999999999
$ws.Range("N3").Value = 46122.6364
$ws.Range("O3").Value = 29545.8038
$ws.Range("B4").Value = 4453.662
$ws.Range("C4").Value = 5543.459
$ws.Range("E4").Value = 4738.919
$ws.Range("F4").Value = 909377.991
$ws.Range("G4").Value = 766298.265
$ws.Range("H4").Value = 143079.726
$ws.Range("I4").Value = 57813.58
$ws.Range("J4").Value = 85266.14599999999
$ws.Range("K4").Value = 782969.8799999999
$ws.Range("L4").Value = 774147.678
$ws.Range("M4").Value = 8822.201999999999
$ws.Range("N4").Value = 5056.351
$ws.Range("O4").Value = 3765.851
